$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.096772333333333
$ws.Range("H2").Value = 3.290317
$ws.Range("I2").Value = 0.2426185621302128
$ws.Range("J2").Value = 0.2426185621302128
$ws.Range("M2").Value = 5.209944
$ws.Range("N2").Value = 15.629832
$ws.Range("O2").Value = 0.1751928672265232
$ws.Range("P2").Value = 0.1751928672265232
$ws.Range("Q2").Value = 5.714122437416
$ws.Range("R2").Value = 51.427101936744
$ws.Range("S2").Value = 0.04250504154196835
$ws.Range("T2").Value = 0.04250504154196835
$ws.Range("G3").Value = 1.096772333333333
$ws.Range("H3").Value = 3.290317
$ws.Range("I3").Value = 0.2426185621302128
$ws.Range("J3").Value = 0.2426185621302128
$ws.Range("M3").Value = 6.497702
$ws.Range("O3").Value = 0.2184958310038485
$ws.Range("P3").Value = 0.2184958310038485
$ws.Range("Q3").Value = 7.126499783844666
$ws.Range("R3").Value = 64.138498054602
$ws.Range("S3").Value = 0.05301114434959969
$ws.Range("T3").Value = 0.05301114434959969
$ws.Range("G4").Value = 1.096772333333333
$ws.Range("H4").Value = 3.290317
$ws.Range("I4").Value = 0.2426185621302128
$ws.Range("J4").Value = 0.2426185621302128
$ws.Range("M4").Value = 1.714656666666667
$ws.Range("N4").Value = 5.14397
$ws.Range("O4").Value = 0.057658127945791
$ws.Range("P4").Value = 0.057658127945791
$ws.Range("Q4").Value = 1.880587993165555
$ws.Range("R4").Value = 16.92529193849
$ws.Range("S4").Value = 0.01398893209732766
$ws.Range("T4").Value = 0.01398893209732766
$ws.Range("G5").Value = 1.096772333333333
$ws.Range("H5").Value = 3.290317
$ws.Range("I5").Value = 0.2426185621302128
$ws.Range("J5").Value = 0.2426185621302128
$ws.Range("M5").Value = 16.31603133333333
$ws.Range("N5").Value = 48.948094
$ws.Range("O5").Value = 0.5486531738238374
$ws.Range("P5").Value = 0.5486531738238374
$ws.Range("Q5").Value = 17.89497175619977
$ws.Range("R5").Value = 161.054745805798
$ws.Range("S5").Value = 0.1331134441413171
$ws.Range("T5").Value = 0.1331134441413171
$ws.Range("I6").Value = 0.03766810132102297
$ws.Range("J6").Value = 0.03766810132102297
$ws.Range("M6").Value = 5.209944
$ws.Range("N6").Value = 15.629832
$ws.Range("O6").Value = 0.1751928672265232
$ws.Range("P6").Value = 0.1751928672265232
$ws.Range("Q6").Value = 0.8871544742640001
$ws.Range("R6").Value = 7.984390268376001
$ws.Range("S6").Value = 0.006599182673409199
$ws.Range("T6").Value = 0.006599182673409199
$ws.Range("I7").Value = 0.03766810132102297
$ws.Range("J7").Value = 0.03766810132102297
$ws.Range("M7").Value = 6.497702
$ws.Range("O7").Value = 0.2184958310038485
$ws.Range("P7").Value = 0.2184958310038485
$ws.Range("R7").Value = 9.957916748358002
$ws.Range("S7").Value = 0.008230323100474075
$ws.Range("T7").Value = 0.008230323100474075
$ws.Range("I8").Value = 0.03766810132102297
$ws.Range("J8").Value = 0.03766810132102297
$ws.Range("M8").Value = 1.714656666666667
$ws.Range("N8").Value = 5.14397
$ws.Range("O8").Value = 0.057658127945791
$ws.Range("P8").Value = 0.057658127945791
$ws.Range("Q8").Value = 0.2919734518566667
$ws.Range("R8").Value = 2.627761066710001
$ws.Range("S8").Value = 0.002171872205442561
$ws.Range("T8").Value = 0.002171872205442561
$ws.Range("I9").Value = 0.03766810132102297
$ws.Range("J9").Value = 0.03766810132102297
$ws.Range("M9").Value = 16.31603133333333
$ws.Range("N9").Value = 48.948094
$ws.Range("O9").Value = 0.5486531738238374
$ws.Range("P9").Value = 0.5486531738238374
$ws.Range("Q9").Value = 2.778310131471333
$ws.Range("R9").Value = 25.004791183242
$ws.Range("S9").Value = 0.02066672334169713
$ws.Range("T9").Value = 0.02066672334169713
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5018676666666667
$ws.Range("H10").Value = 1.505603
$ws.Range("I10").Value = 0.1110188577571507
$ws.Range("J10").Value = 0.1110188577571507
$ws.Range("M10").Value = 5.209944
$ws.Range("N10").Value = 15.629832
$ws.Range("O10").Value = 0.1751928672265232
$ws.Range("P10").Value = 0.1751928672265232
$ws.Range("Q10").Value = 2.614702438744
$ws.Range("R10").Value = 23.532321948696
$ws.Range("S10").Value = 0.01944971200668877
$ws.Range("T10").Value = 0.01944971200668877
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5018676666666667
$ws.Range("H11").Value = 1.505603
$ws.Range("I11").Value = 0.1110188577571507
$ws.Range("J11").Value = 0.1110188577571507
$ws.Range("M11").Value = 6.497702
$ws.Range("O11").Value = 0.2184958310038485
$ws.Range("P11").Value = 0.2184958310038485
$ws.Range("Q11").Value = 3.260986541435334
$ws.Range("R11").Value = 29.348878872918
$ws.Range("S11").Value = 0.0242571575827467
$ws.Range("T11").Value = 0.02425715758274669
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5018676666666667
$ws.Range("H12").Value = 1.505603
$ws.Range("I12").Value = 0.1110188577571507
$ws.Range("J12").Value = 0.1110188577571507
$ws.Range("M12").Value = 1.714656666666667
$ws.Range("N12").Value = 5.14397
$ws.Range("O12").Value = 0.057658127945791
$ws.Range("P12").Value = 0.057658127945791
$ws.Range("Q12").Value = 0.8605307404344446
$ws.Range("R12").Value = 7.744776663910001
$ws.Range("S12").Value = 0.006401139504957369
$ws.Range("T12").Value = 0.006401139504957368
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5018676666666667
$ws.Range("H13").Value = 1.505603
$ws.Range("I13").Value = 0.1110188577571507
$ws.Range("J13").Value = 0.1110188577571507
$ws.Range("M13").Value = 16.31603133333333
$ws.Range("N13").Value = 48.948094
$ws.Range("O13").Value = 0.5486531738238374
$ws.Range("P13").Value = 0.5486531738238374
$ws.Range("Q13").Value = 8.188488574520221
$ws.Range("R13").Value = 73.69639717068199
$ws.Range("S13").Value = 0.06091084866275789
$ws.Range("T13").Value = 0.06091084866275788
$ws.Range("G14").Value = 2.751641333333334
$ws.Range("H14").Value = 8.254924000000001
$ws.Range("I14").Value = 0.6086944787916135
$ws.Range("J14").Value = 0.6086944787916135
$ws.Range("M14").Value = 5.209944
$ws.Range("N14").Value = 15.629832
$ws.Range("O14").Value = 0.1751928672265232
$ws.Range("P14").Value = 0.1751928672265232
$ws.Range("Q14").Value = 14.335897254752
$ws.Range("R14").Value = 129.023075292768
$ws.Range("S14").Value = 0.1066389310044569
$ws.Range("T14").Value = 0.1066389310044569
$ws.Range("G15").Value = 2.751641333333334
$ws.Range("H15").Value = 8.254924000000001
$ws.Range("I15").Value = 0.6086944787916135
$ws.Range("J15").Value = 0.6086944787916135
$ws.Range("M15").Value = 6.497702
$ws.Range("O15").Value = 0.2184958310038485
$ws.Range("P15").Value = 0.2184958310038485
$ws.Range("Q15").Value = 17.87934539488267
$ws.Range("R15").Value = 160.914108553944
$ws.Range("S15").Value = 0.132997205971028
$ws.Range("T15").Value = 0.132997205971028
$ws.Range("G16").Value = 2.751641333333334
$ws.Range("H16").Value = 8.254924000000001
$ws.Range("I16").Value = 0.6086944787916135
$ws.Range("J16").Value = 0.6086944787916135
$ws.Range("M16").Value = 1.714656666666667
$ws.Range("N16").Value = 5.14397
$ws.Range("O16").Value = 0.057658127945791
$ws.Range("P16").Value = 0.057658127945791
$ws.Range("Q16").Value = 4.718120156475556
$ws.Range("R16").Value = 42.46308140828
$ws.Range("S16").Value = 0.03509618413806342
$ws.Range("T16").Value = 0.03509618413806342
$ws.Range("G17").Value = 2.751641333333334
$ws.Range("H17").Value = 8.254924000000001
$ws.Range("I17").Value = 0.6086944787916135
$ws.Range("J17").Value = 0.6086944787916135
$ws.Range("M17").Value = 16.31603133333333
$ws.Range("N17").Value = 48.948094
$ws.Range("O17").Value = 0.5486531738238374
$ws.Range("P17").Value = 0.5486531738238374
$ws.Range("Q17").Value = 44.89586621276177
$ws.Range("R17").Value = 404.062795914856
$ws.Range("S17").Value = 0.3339621576780653
$ws.Range("T17").Value = 0.3339621576780653
